$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Email column header and value
$ws.Range("G1").Value = "Email"
$ws.Range("G2").Value = "adrianrentea01@gmail.com"

# Set column width for new column G (stored width ends up 24.5, matching
# the other best-fit columns already on the sheet)
$ws.Columns.Item(7).ColumnWidth = 23.6666666666667

# Update selection to match target state
$ws.Range("G1:G2").Select()
